$d = $word.ActiveDocument
$d.Content.Find.Execute("normal PyTorch.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "normal PyTorch. The graph below shows CPU Wall time comparison between Pytorch vs IPEX.", 2)
